# Update cryptocurrency price/volume snapshot values in the worksheet.
# (GitHub Actions scheduled data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.528.42'
$ws.Range("E2").Value = '  +1.06%  '

$ws.Range("D3").Value = '2.485.60'
$ws.Range("E3").Value = '  +1.15%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.23'
$ws.Range("E5").Value = '  +0.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.15'
$ws.Range("E6").Value = '  -0.95%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.545'
$ws.Range("E7").Value = '  -0.75%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.497'
$ws.Range("E9").Value = '  -0.80%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.63'
$ws.Range("E10").Value = '  -2.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0781'
$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("D13").Value = '2.866.23'
$ws.Range("E13").Value = '  +1.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.84'
$ws.Range("E14").Value = '  -1.39%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.507.79'
$ws.Range("E15").Value = '  +2.21%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.43'
$ws.Range("E16").Value = '  +6.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.753'
$ws.Range("E17").Value = '  -4.02%  '

$ws.Range("D18").Value = '41.654.10'
$ws.Range("E18").Value = '  +1.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.33'
$ws.Range("E19").Value = '  -0.37%  '

$ws.Range("D20").Value = '0.0₃0921'
$ws.Range("E20").Value = '  +0.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.58'
$ws.Range("E21").Value = '  +5.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.14'
$ws.Range("E22").Value = '  -3.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.45'
$ws.Range("E23").Value = '  -0.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.70'
$ws.Range("E24").Value = '  -2.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("E26").Value = '  -1.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.91'
$ws.Range("E27").Value = '  +1.81%  '

$ws.Range("E28").Value = '  +0.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.64'
$ws.Range("E29").Value = '  -0.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.21'
$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '154.98'
$ws.Range("E31").Value = '  +1.38%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.41'
$ws.Range("E32").Value = '  -2.87%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.57'
$ws.Range("E33").Value = '  -0.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.12'
$ws.Range("E34").Value = '  +6.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0754'
$ws.Range("E35").Value = '  +0.45%  '

$ws.Range("E36").Value = '  -3.04%  '

$ws.Range("E37").Value = '  -1.93%  '

$ws.Range("E38").Value = '  -3.52%  '

$ws.Range("E39").Value = '  +0.13%  '

$ws.Range("E41").Value = '  -1.81%  '

$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.95'
$ws.Range("E43").Value = '  -5.30%  '

$ws.Range("D44").Value = '1.960.84'
$ws.Range("E44").Value = '  -0.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0284'
$ws.Range("E45").Value = '  -0.11%  '

$ws.Range("E46").Value = '  -3.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.87'
$ws.Range("E47").Value = '  +1.49%  '

$ws.Range("D48").Value = '2.725.19'
$ws.Range("E48").Value = '  +1.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '96.16'
$ws.Range("E49").Value = '  -1.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.27'
$ws.Range("E50").Value = '  -3.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.24'
$ws.Range("E51").Value = '  -3.96%  '
